# Update odds values on the active worksheet to match the new FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.88
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 1.48
$ws.Range("U2").Value = 4.18
$ws.Range("V2").Value = 1.21
$ws.Range("AA2").Value = 2.1
$ws.Range("AB2").Value = 1.67

# Row 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85
$ws.Range("U3").Value = 4.12

# Row 4
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.3
$ws.Range("L4").Value = 6.5
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("U4").Value = 3.8
$ws.Range("V4").Value = 1.26

# Row 5
$ws.Range("G5").Value = 1.91
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.63
$ws.Range("L5").Value = 4.75
$ws.Range("W5").Value = 4
$ws.Range("X5").Value = 1.22
$ws.Range("AA5").Value = 2
$ws.Range("AB5").Value = 1.73

# Row 7
$ws.Range("G7").Value = 2.8
$ws.Range("H7").Value = 2.85
$ws.Range("I7").Value = 2.62
$ws.Range("J7").Value = 3.35
$ws.Range("K7").Value = 1.93
$ws.Range("L7").Value = 3.3
$ws.Range("W7").Value = 3.6
$ws.Range("Y7").Value = 1.47
$ws.Range("Z7").Value = 2.35
